$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change B2:B6 values from numeric 123 to the text string "abc"
$ws.Range("B2:B6").Value = "abc"

# Update selection to B3:B6 with active cell B3
$ws.Range("B3:B6").Select()
